$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.975.65"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "2.232.31"
$ws.Range("E3").Value = "  +3.98%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'252.34"
$ws.Range("E5").Value = "  +6.80%  "
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").Value = "'75.06"
$ws.Range("E7").Value = "  +7.20%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").Value = "'41.26"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").Value = "'0.0924"
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("E12").Value = "  +3.99%  "
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").Value = "2.567.88"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "'14.51"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "2.229.63"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "42.885.17"
$ws.Range("E18").Value = "  +4.42%  "
$ws.Range("E19").Value = "  +3.29%  "
$ws.Range("D20").Value = "'71.22"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("E21").Value = "  +4.06%  "
$ws.Range("D22").Value = "'229.85"
$ws.Range("E22").Value = "  +2.06%  "
$ws.Range("E23").Value = "  +12.53%  "
$ws.Range("D24").Value = "'9.34"
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'10.71"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "'39.03"
$ws.Range("E28").Value = "  +24.25%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("D31").Value = "'170.22"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +2.67%  "
$ws.Range("D33").Value = "'0.0795"
$ws.Range("E33").Value = "  +5.27%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("E35").Value = "  +11.86%  "
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "'4.46"
$ws.Range("E37").Value = "  +5.45%  "
$ws.Range("E38").Value = "  +12.07%  "
$ws.Range("D39").Value = "'12.37"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").Value = "'2.11"
$ws.Range("E40").Value = "  +2.96%  "
$ws.Range("E41").Value = "  +9.35%  "
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "'59.69"
$ws.Range("E43").Value = "  +3.18%  "
$ws.Range("E44").Value = "  +28.57%  "
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("D46").Value = "'103.21"
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("D47").Value = "'0.0984"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +13.35%  "
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("E51").Value = "  +2.35%  "

# Reset number-formatted cells back to the default (unstyled) appearance
# now that the text has been forced in, so no stray style survives.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
